$wb = $excel.ActiveWorkbook

# Update the Company sheet value from "International Wire Group, Inc." to "StandardTestCompany"
$wsCompany = $wb.Worksheets.Item("Company")
$wsCompany.Range("A2").Value = "StandardTestCompany"

# Update selection on ActivityColumns sheet (keep A5 selected) before switching away
$wsActivityColumns = $wb.Worksheets.Item("ActivityColumns")
$wsActivityColumns.Activate()
$wsActivityColumns.Range("A5").Select()

# Make Company the active sheet and select F8
$wsCompany.Activate()
$wsCompany.Range("F8").Select()
